$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G2").Value = 300
$ws1.Range("D12").Select()

# --- Sheet3 ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C2").Formula = "=890"
$ws3.Range("E2").Value = 2
$ws3.Range("E3").Value = 2
$ws3.Range("E4").Value = 2
$ws3.Range("E14").Select()

# --- Sheet4 ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("C2").Formula = "=890"
$ws4.Range("E2").Value = 2
$ws4.Range("E3").Value = 2
$ws4.Range("E4").Value = 2
$ws4.Range("G17").Select()
